# Apply the changes described by the diff to test7.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on C2 (mailto:juan@example.com)
ForEach ($hl in $ws.Hyperlinks) {
    $hl.Delete()
}

# Drop the now-unused trailing columns F:I (the sheet shrinks from A1:I2 to A1:E2)
$ws.Range("F1:I2").Clear()

# Rewrite the header row (row 1)
$ws.Range("A1").Value = "Nombre"
$ws.Range("B1").Value = "locacalizacion"
$ws.Range("C1").Value = "Correo electrónico"
$ws.Range("D1").Value = "id"
$ws.Range("E1").Value = "kind"

# Rewrite the data row (row 2)
$ws.Range("A2").Value = "jorge"
$ws.Range("B2").Value = "18:13:14:12S"
$ws.Range("C2").Value = "jorge@email.es"
$ws.Range("D2").Value = "ID4"
$ws.Range("E2").Value = 1

# Update the selected range to match the new used range
$ws.Range("A1:E2").Select() | Out-Null
